$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- First, extend sheet1 ("1") with the new column G and new rows 5-11 ---

# New header cell G1
$ws1.Range("G1").Value = "Katze"

# New G column values for existing rows 2-4
$ws1.Range("G2").Formula = '="0080"'
$ws1.Range("G3").Value = 1
$ws1.Range("G4").Value = 4

# New rows 5-11: copy formatting from row 4 (B4:G4) down, then fill in data
$ws1.Range("B4:G4").Copy()
$ws1.Range("B5:G11").PasteSpecial(-4122)

# Row 5
$ws1.Range("B5").Formula = '="0300"'
$ws1.Range("C5").Formula = '="1.2"'
$ws1.Range("D5").Value = "Some Other Item"
$ws1.Range("E5").Value = 67
$ws1.Range("F5").Value = 412
$ws1.Range("G5").Value = 6

# Row 6
$ws1.Range("B6").Formula = '="0400"'
$ws1.Range("D6").Value = "Some Other Item"
$ws1.Range("E6").Value = 23
$ws1.Range("F6").Value = 65
$ws1.Range("G6").Value = 1

# Row 7
$ws1.Range("B7").Formula = '="0500"'
$ws1.Range("D7").Value = "Some Other Item"
$ws1.Range("E7").Value = 6
$ws1.Range("F7").Value = 3
$ws1.Range("G7").Value = 2

# Row 8
$ws1.Range("B8").Formula = '="0600"'
$ws1.Range("D8").Value = "Some Other Item"
$ws1.Range("E8").Value = 1
$ws1.Range("F8").Value = 1
$ws1.Range("G8").Value = 6

# Row 9
$ws1.Range("B9").Formula = '="0700"'
$ws1.Range("D9").Value = "Some Other Item"
$ws1.Range("E9").Value = 64
$ws1.Range("F9").Value = 66
$ws1.Range("G9").Value = 7

# Row 10
$ws1.Range("B10").Formula = '="0800"'
$ws1.Range("D10").Value = "Some Other Item"
$ws1.Range("E10").Value = 3
$ws1.Range("F10").Value = 345
$ws1.Range("G10").Value = 2

# Row 11
$ws1.Range("B11").Formula = '="0900"'
$ws1.Range("D11").Value = "Some Other Item"
$ws1.Range("E11").Value = 8
$ws1.Range("F11").Value = 12
$ws1.Range("G11").Value = 33

# C5:C11 share the same "1.2" formula -> set together so engine stores it as one shared formula
$ws1.Range("C5:C11").Formula = '="1.2"'

# --- Rename sheet1 from "1" to "Eins" ---
$ws1.Name = "Eins"

# --- Duplicate the sheet to create "Zwei" ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Zwei"

# Sheet "Zwei" differs from "Eins" only in that E4 is empty
$ws2.Range("E4").ClearContents()

# --- Selections ---
$ws2.Range("F7").Select()
$ws1.Activate()
$ws1.Range("D8").Select()
